# Update mark for HaiTCT
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Report names first (A column), in the order the new strings were authored
$ws.Range("A45").Value = "Report 01 (Review + modify)"
$ws.Range("A48").Value = "Report 04 ( Page 42-46, 56-74)"
$ws.Range("A47").Value = "Report 03 SRS_PMS (Page 19-32, 39-64), User Requirement_PMS(4.Project Eye, 7.Admin)"
$ws.Range("A46").Value = "Report 02 (Review)"

# Note column (C), all the same new string "Finish task on time"
$ws.Range("C45").Value = "Finish task on time"
$ws.Range("C46").Value = "Finish task on time"
$ws.Range("C47").Value = "Finish task on time"
$ws.Range("C48").Value = "Finish task on time"

# Mark column (B), numeric value 10 for all four report rows
$ws.Range("B45").Value = 10
$ws.Range("B46").Value = 10
$ws.Range("B47").Value = 10
$ws.Range("B48").Value = 10

# Left-align the section label cell for the HaiTCT block header
$ws.Range("A44").HorizontalAlignment = -4131

# Update the view selection to match where the editor was last working
$ws.Activate()
$ws.Range("C41").Select()
